# Apply the "Cookie Enhance Lv, Star api" change:
# Insert a new "POINT_COOKIE_LV" entry (value 102) into the EObjType table
# on the "ObjType" sheet, right before "POINT_C_GACHA_NORMAL", and make the
# "ObjType" sheet the active tab/selection instead of "Common".

$wb = $excel.ActiveWorkbook
$wsCommon  = $wb.Worksheets.Item("Common")
$wsObjType = $wb.Worksheets.Item("ObjType")

# Insert a new row above row 10 (POINT_C_GACHA_NORMAL) to make room for the
# new enum entry, shifting everything below it down by one row.
$wsObjType.Rows.Item(10).Insert()

# Fill in the new row's data.
$wsObjType.Range("A10").Value = "EObjType"
$wsObjType.Range("B10").Value = "POINT_COOKIE_LV"
$wsObjType.Range("C10").Value = 102

# Update selection / active state: ObjType tab becomes active, Common is
# no longer the selected tab.
$wsObjType.Range("C11").Select()
$wsObjType.Activate()

$wb.Save()
